# Daily attendance processing - 2025-11-03 18:51:15
#
# Column G ("Recorded By") lists who/what recorded each attendance
# session as a comma-separated string. Re-order each such list so the
# "System" entry is listed first, leaving any list that does not
# include "System" (or that is a single value) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -ne $null -and $v -ne "") {
        $parts = $v -split ', '

        if ($parts.Length -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p -eq 'System') {
                    $hasSystem = $true
                }
            }

            if ($hasSystem) {
                $reversed = $parts[($parts.Length - 1)..0]
                $cell.Value2 = $reversed -join ', '
            }
        }
    }
}
